# Update order and authentication
# Replace the menu item rows (2-9) with the new crawled data, leaving the
# header row (1) and the "extra items" rows (10-14) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Mỳ quảng chay"
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr361idjlk7t4b"

# Row 3
$ws.Range("A3").Value = "Phile cá chiên sốt teriyaki"
$ws.Range("B3").Value = "• Cơm gạo dẻo, hộp 4 ngăn tiện lợi. • Giá đã bao gồm đủ đồ xào và canh ăn kèm. • Rất vui được ăn trưa cùng bạn."
$ws.Range("C3").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lqy3g01uu6ys91"

# Row 4
$ws.Range("A4").Value = "Trứng chiên thịt bằm"
$ws.Range("B4").Value = "• Cơm gạo dẻo, nhiều rau xanh`n• Có canh ăn kèm miễn phí`n• Hộp đựng 4 ngăn tiện lợi"
$ws.Range("C4").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lqynexwnwz44e1"

# Row 5
$ws.Range("A5").Value = "Gà xào bắp non nấm bào ngư"
$ws.Range("B5").Value = "Cơm gà chiên ăn kèm cơm, đồ xào, canh"
$ws.Range("C5").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr30t0ay6z0k75"

# Row 6
$ws.Range("A6").Value = "Sườn cốt lết nướng mật ong"
$ws.Range("B6").Value = "• Cơm gạo dẻo, nhiều rau xanh. • Có canh ăn kèm miễn phí. • Hộp đựng 4 ngăn tiện dùng."
$ws.Range("C6").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr5ynxhsng0418"

# Row 7
$ws.Range("A7").Value = "Gà om ớt hiêm"
$ws.Range("B7").Value = "Cơm + Đồ mặn + Đồ xào + Canh + Trái cây"
$ws.Range("C7").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr320g6zmojo64"

# Row 8
$ws.Range("A8").Value = "Canh chua cá lóc"
$ws.Range("B8").Value = $null
$ws.Range("C8").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr343zj47jyc01"

# Row 9
$ws.Range("A9").Value = "Nui sốt bò bằm"
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr2w1utsrnpg99"

$wb.Save()
